$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.547311067581177
$ws.Range("B1").Value = 2.377239942550659
$ws.Range("C1").Value = 1.822073459625244
$ws.Range("D1").Value = 1.647401332855225
$ws.Range("E1").Value = 1.484955549240112
